$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 172.75
$ws.Range("I8").Value = 172.75
$ws.Range("K8").Value = 518.25
$ws.Range("M8").Value = -379.25

$ws.Range("H17").Value = 513990.16
$ws.Range("J17").Value = 513990.16
$ws.Range("L17").Value = 1541970.48
$ws.Range("N17").Value = -1542306.48

$ws.Range("H19").Value = 1223.9048
$ws.Range("I19").Value = 1530.9231
$ws.Range("K19").Value = 1530.9231
$ws.Range("M19").Value = -1355.9231

$ws.Range("H21").Value = 4401.8887
$ws.Range("I21").Value = 2323.4
$ws.Range("K21").Value = 2323.4
$ws.Range("M21").Value = -1855.4

$ws.Range("H23").Value = 4401.8887
$ws.Range("I23").Value = 2323.4
$ws.Range("K23").Value = 2323.4
$ws.Range("M23").Value = -2089.4

$ws.Range("H38").Value = 24.875
$ws.Range("I38").Value = 29
$ws.Range("J38").Value = 18
$ws.Range("K38").Value = 87
$ws.Range("L38").Value = 54
$ws.Range("M38").Value = 285
$ws.Range("N38").Value = -798

$ws.Range("H42").Value = 78.5
$ws.Range("I42").Value = 51.857143
$ws.Range("J42").Value = 105.14286
$ws.Range("K42").Value = 155.571429
$ws.Range("L42").Value = 315.42858
$ws.Range("M42").Value = 74.42857100000001
$ws.Range("N42").Value = -775.42858

$ws.Range("H103").Value = 886.125
$ws.Range("I103").Value = 888.4545000000001
$ws.Range("J103").Value = 881
$ws.Range("K103").Value = 2665.3635
$ws.Range("L103").Value = 2643
$ws.Range("M103").Value = -2079.3635
$ws.Range("N103").Value = -3815

$ws.Range("H138").Value = 2539.64
$ws.Range("I138").Value = 1138.6129
$ws.Range("J138").Value = 3169.087
$ws.Range("K138").Value = 3415.8387
$ws.Range("L138").Value = 9507.261
$ws.Range("M138").Value = 1724.1613
$ws.Range("N138").Value = -19787.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1672.9048
$ws.Range("I122").Value = 1369.5
$ws.Range("J122").Value = 2077.4443
$ws.Range("K122").Value = 4108.5
$ws.Range("L122").Value = 6232.3329
$ws.Range("M122").Value = -1658.5
$ws.Range("N122").Value = -11132.3329

$ws.Range("H133").Value = 45000
$ws.Range("J133").Value = 45000
$ws.Range("L133").Value = 45000
$ws.Range("N133").Value = -50060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1620.7059
$ws.Range("I94").Value = 1675.8572
$ws.Range("J94").Value = 1363.3334
$ws.Range("K94").Value = 1675.8572
$ws.Range("L94").Value = 1363.3334
$ws.Range("M94").Value = -1224.8572
$ws.Range("N94").Value = -2265.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 301.5
$ws.Range("I7").Value = 304.125
$ws.Range("J7").Value = 291
$ws.Range("K7").Value = 304.125
$ws.Range("L7").Value = 291
$ws.Range("M7").Value = -191.125
$ws.Range("N7").Value = -517

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = ""

$ws.Range("H22").Value = 1047.5
$ws.Range("J22").Value = 248.33333
$ws.Range("L22").Value = 248.33333
$ws.Range("N22").Value = -948.3333299999999

$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").Value = ""

$ws.Range("H31").Value = 2091.2678
$ws.Range("I31").Value = 2213.15
$ws.Range("J31").Value = 1786.5625
$ws.Range("K31").Value = 2213.15
$ws.Range("L31").Value = 1786.5625
$ws.Range("M31").Value = -1918.15
$ws.Range("N31").Value = -2376.5625

$ws.Range("H34").Value = 2091.2678
$ws.Range("I34").Value = 2213.15
$ws.Range("J34").Value = 1786.5625
$ws.Range("K34").Value = 2213.15
$ws.Range("L34").Value = 1786.5625
$ws.Range("M34").Value = -2011.15
$ws.Range("N34").Value = -2190.5625

$ws.Range("H122").Value = 1710
$ws.Range("I122").Value = 933.3333
$ws.Range("J122").Value = 2875
$ws.Range("K122").Value = 2799.9999
$ws.Range("L122").Value = 8625
$ws.Range("M122").Value = -349.9998999999998
$ws.Range("N122").Value = -13525

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = ""

$ws.Range("H132").Value = 1711.6364
$ws.Range("I132").Value = 860
$ws.Range("J132").Value = 2421.3333
$ws.Range("K132").Value = 2580
$ws.Range("L132").Value = 7263.999899999999
$ws.Range("M132").Value = -50
$ws.Range("N132").Value = -12323.9999

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = ""

$ws.Range("H137").Value = 79800
$ws.Range("J137").Value = 79800
$ws.Range("L137").Value = 79800
$ws.Range("N137").Value = -90000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8081710.5
$ws.Range("I131").Value = 71573096
$ws.Range("J131").Value = 989.16364
$ws.Range("K131").Value = 214719288
$ws.Range("L131").Value = 2967.49092
$ws.Range("M131").Value = -214714248
$ws.Range("N131").Value = -13047.49092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 65004.668
$ws.Range("I23").Value = 30000
$ws.Range("J23").Value = 82507
$ws.Range("K23").Value = 30000
$ws.Range("L23").Value = 82507
$ws.Range("M23").Value = -29770
$ws.Range("N23").Value = -82967

$ws.Range("H115").Value = 14651
$ws.Range("J115").Value = 14651
$ws.Range("L115").Value = 14651
$ws.Range("N115").Value = -17001

$ws.Range("H121").Value = 11863.728
$ws.Range("J121").Value = 11863.728
$ws.Range("L121").Value = 11863.728
$ws.Range("N121").Value = -15357.728

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1378.0555
$ws.Range("I122").Value = 1236.0714
$ws.Range("J122").Value = 1875
$ws.Range("K122").Value = 3708.2142
$ws.Range("L122").Value = 5625
$ws.Range("M122").Value = -1258.2142
$ws.Range("N122").Value = -10525

$ws.Range("H126").Value = 425.26315
$ws.Range("I126").Value = 452.94116
$ws.Range("K126").Value = 1358.82348
$ws.Range("M126").Value = 1111.17652
